$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2281.7
$ws.Range("I15").Value = 2281.7
$ws.Range("K15").Value = 6845.099999999999
$ws.Range("M15").Value = -6676.099999999999
$ws.Range("H41").Value = 975
$ws.Range("I41").Value = 975
$ws.Range("K41").Value = 975
$ws.Range("M41").Value = -535
$ws.Range("H70").Value = 3712.5
$ws.Range("I70").Value = 3406.25
$ws.Range("K70").Value = 10218.75
$ws.Range("M70").Value = -9948.75
$ws.Range("H73").Value = 3712.5
$ws.Range("I73").Value = 3406.25
$ws.Range("K73").Value = 10218.75
$ws.Range("M73").Value = -9282.75
$ws.Range("H80").Value = 208.66667
$ws.Range("I80").Value = 228.5
$ws.Range("J80").Value = 50
$ws.Range("K80").Value = 685.5
$ws.Range("L80").Value = 150
$ws.Range("M80").Value = 312.5
$ws.Range("N80").Value = -2146
$ws.Range("H83").Value = 208.66667
$ws.Range("I83").Value = 228.5
$ws.Range("J83").Value = 50
$ws.Range("K83").Value = 2056.5
$ws.Range("L83").Value = 450
$ws.Range("M83").Value = 2935.5
$ws.Range("N83").Value = -10434
$ws.Range("H92").Value = 45454950
$ws.Range("I92").Value = 62500396
$ws.Range("J92").Value = 433.5
$ws.Range("K92").Value = 62500396
$ws.Range("L92").Value = 433.5
$ws.Range("M92").Value = -62499148
$ws.Range("N92").Value = -2929.5
$ws.Range("H95").Value = 25686.25
$ws.Range("J95").Value = 25686.25
$ws.Range("L95").Value = 25686.25
$ws.Range("N95").Value = -31178.25
$ws.Range("H107").Value = 61120.8
$ws.Range("I107").Value = 65425.855
$ws.Range("K107").Value = 65425.855
$ws.Range("M107").Value = -63505.855
$ws.Range("H113").Value = 4065
$ws.Range("I113").Value = 3997.5
$ws.Range("J113").Value = 4200
$ws.Range("K113").Value = 3997.5
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = -743.5
$ws.Range("N113").Value = -10708
$ws.Range("H127").Value = 1048.25
$ws.Range("I127").Value = 1048.25
$ws.Range("K127").Value = 3144.75
$ws.Range("M127").Value = 1815.25
$ws.Range("H137").Value = 2035.7142
$ws.Range("I137").Value = 625
$ws.Range("K137").Value = 1875
$ws.Range("M137").Value = 675
$ws.Range("H138").Value = 3770.4736
$ws.Range("J138").Value = 4265
$ws.Range("L138").Value = 12795
$ws.Range("N138").Value = -23075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.5
$ws.Range("I45").Value = 1999
$ws.Range("K45").Value = 1999
$ws.Range("M45").Value = -1622
$ws.Range("H63").Value = 1995
$ws.Range("I63").Value = 1995
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1995
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1309
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1995
$ws.Range("I66").Value = 1995
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9975
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6543
$ws.Range("N66").ClearContents()
$ws.Range("H102").Value = 26250938
$ws.Range("I102").Value = 1429644.1
$ws.Range("K102").Value = 1429644.1
$ws.Range("M102").Value = -1428022.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 13499.5
$ws.Range("J95").Value = 13499.5
$ws.Range("L95").Value = 13499.5
$ws.Range("N95").Value = -18991.5
$ws.Range("H107").Value = 100815.75
$ws.Range("I107").Value = 200387.5
$ws.Range("K107").Value = 200387.5
$ws.Range("M107").Value = -198467.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 739.7143
$ws.Range("J16").Value = 766.3333
$ws.Range("L16").Value = 766.3333
$ws.Range("N16").Value = -1340.3333
$ws.Range("H41").Value = 3666.3333
$ws.Range("I41").Value = 3666.3333
$ws.Range("K41").Value = 3666.3333
$ws.Range("M41").Value = -3238.3333
$ws.Range("H47").Value = 32714
$ws.Range("J47").Value = 31499.666
$ws.Range("L47").Value = 31499.666
$ws.Range("N47").Value = -32631.666
$ws.Range("H99").Value = 387004.84
$ws.Range("I99").Value = 558038.75
$ws.Range("J99").Value = 2178.5
$ws.Range("K99").Value = 558038.75
$ws.Range("L99").Value = 2178.5
$ws.Range("M99").Value = -556540.75
$ws.Range("N99").Value = -5174.5
$ws.Range("H107").Value = 837.1818
$ws.Range("I107").Value = 701.5714
$ws.Range("K107").Value = 701.5714
$ws.Range("M107").Value = 1218.4286
$ws.Range("H113").Value = 739.7143
$ws.Range("J113").Value = 766.3333
$ws.Range("L113").Value = 766.3333
$ws.Range("N113").Value = -5106.3333
$ws.Range("H126").Value = 387004.84
$ws.Range("I126").Value = 558038.75
$ws.Range("J126").Value = 2178.5
$ws.Range("K126").Value = 1674116.25
$ws.Range("L126").Value = 6535.5
$ws.Range("M126").Value = -1671646.25
$ws.Range("N126").Value = -11475.5
$ws.Range("H141").Value = 787212.5
$ws.Range("J141").Value = 787212.5
$ws.Range("L141").Value = 787212.5
$ws.Range("N141").Value = -797572.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H26").Value = 372.375
$ws.Range("J26").Value = 990
$ws.Range("L26").Value = 2970
$ws.Range("N26").Value = -3546
$ws.Range("H32").Value = 8000
$ws.Range("J32").Value = 8000
$ws.Range("L32").Value = 24000
$ws.Range("N32").Value = -24566
$ws.Range("H44").Value = 816.3889
$ws.Range("I44").Value = 198.2
$ws.Range("K44").Value = 594.5999999999999
$ws.Range("M44").Value = -196.5999999999999
$ws.Range("H55").Value = 3902.3572
$ws.Range("J55").Value = 3953.4783
$ws.Range("L55").Value = 11860.4349
$ws.Range("N55").Value = -12214.4349
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H94").Value = 450
$ws.Range("I94").Value = 450
$ws.Range("K94").Value = 1350
$ws.Range("M94").Value = -674
$ws.Range("H117").Value = 2998.4614
$ws.Range("I117").Value = 423.33334
$ws.Range("J117").Value = 3771
$ws.Range("K117").Value = 1270.00002
$ws.Range("L117").Value = 11313
$ws.Range("M117").Value = 2171.99998
$ws.Range("N117").Value = -18197

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 36092.8
$ws.Range("I55").Value = 35199.332
$ws.Range("J55").Value = 37433
$ws.Range("K55").Value = 35199.332
$ws.Range("L55").Value = 37433
$ws.Range("M55").Value = -34872.332
$ws.Range("N55").Value = -38087
$ws.Range("H87").Value = 100000
$ws.Range("I87").Value = 100000
$ws.Range("K87").Value = 100000
$ws.Range("M87").Value = -98752
$ws.Range("H90").Value = 100000
$ws.Range("I90").Value = 100000
$ws.Range("K90").Value = 300000
$ws.Range("M90").Value = -293760
$ws.Range("H94").Value = 61999.668
$ws.Range("I94").Value = 61999
$ws.Range("K94").Value = 61999
$ws.Range("M94").Value = -61323
$ws.Range("H97").Value = 1630.9
$ws.Range("I97").Value = 1474.8667
$ws.Range("K97").Value = 1474.8667
$ws.Range("M97").Value = -978.8667
$ws.Range("H102").Value = 838.2222
$ws.Range("I102").Value = 691.25
$ws.Range("J102").Value = 2014
$ws.Range("K102").Value = 691.25
$ws.Range("L102").Value = 2014
$ws.Range("M102").Value = 930.75
$ws.Range("N102").Value = -5258
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").ClearContents()
$ws.Range("H113").Value = 2454.2727
$ws.Range("I113").Value = 1774.5
$ws.Range("J113").Value = 2842.7144
$ws.Range("K113").Value = 1774.5
$ws.Range("L113").Value = 2842.7144
$ws.Range("M113").Value = 395.5
$ws.Range("N113").Value = -7182.7144
$ws.Range("H122").Value = 1076
$ws.Range("I122").Value = 1090.5714
$ws.Range("K122").Value = 3271.7142
$ws.Range("M122").Value = -821.7142000000003
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 19913.25
$ws.Range("J136").Value = 19913.25
$ws.Range("L136").Value = 59739.75
$ws.Range("N136").Value = -64839.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 29999.5
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -20814
$ws.Range("H46").Value = 3345.7144
$ws.Range("J46").Value = 4187
$ws.Range("L46").Value = 4187
$ws.Range("N46").Value = -4563
$ws.Range("H136").Value = 5400
$ws.Range("J136").Value = 800
$ws.Range("L136").Value = 2400
$ws.Range("N136").Value = -7500

